# Apply the cryptos-list refresh: updated prices / 1h volume deltas,
# plus a couple of rank swaps and one newly-listed coin (BabyDogeCoin)
# that bumped Algorand/Cronos down a row and dropped EnergySwap off the
# bottom of the (fixed 50-row) list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All B/C/D/E cells hold text in this sheet (inline strings, General
# format). Force Text format before writing so Excel does not
# reinterpret numeric-looking strings like "27.945.22" or "0.100" as
# numbers (which would mangle the thousand-dot formatting / drop
# trailing zeros).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.945.22"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.17%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.640.43"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.85"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.524"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.42%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.48"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.29%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.01%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.29%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.875.17"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.653.01"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.62%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.55%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.570"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.54"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.02%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.954.36"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.19"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.42%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.33%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.59"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.98%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.54"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.71%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.37"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.56%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.06"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.90"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.68"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.13%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.24%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.35%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.51%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.40"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.00%  "

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.408.16"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.20%  "

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.09"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.24%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.29%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.81%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.563"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.51%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.928"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.879"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.98%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.91%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.86"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.90%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "67.15"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.52"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.59%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.782.76"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.99"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.18%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0103"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.31%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.14%  "

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0506"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.34%  "
